$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tolerances")

# Set the value of D11 to 70 (MINIMUM accepted average force %)
$ws.Range("D11").Value = 70

# Update the active selection to C17
$ws.Range("C17").Select()
